$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.628.03"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "2.487.82"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.114"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.365"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.98%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000185"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.15%  "
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "63.568.84"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "2.499.32"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.63%  "
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "644.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.35%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000106"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.05%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "2.613.18"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("E29").Value = "  +8.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.995"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.388"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.32%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "152.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0551"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.48%  "
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.748"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.95%  "
